$wb = $excel.ActiveWorkbook

# --- Update the daily conversion summary text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value()
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 7.12 = 28511.6 pesos"), "✅ 1000 Bs = 6.99 = 27993.08 pesos"
$text = $text -replace [regex]::Escape("✅ 28511.6 pesos = 7.09 = 951.84 Bs"), "✅ 27993.08 pesos = 6.97 = 963.84 Bs"
$wsHoja1.Range("A1").Value = $text

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 143
$wsTasas.Range("O10").Value = 4003.01
$wsTasas.Range("N12").Value = 4016.99
$wsTasas.Range("O12").Value = 138.31
